$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 (shifts rows 44-121 down to 45-122),
# adding the "Waypoint" translation entry (object.WAYPOINT / Waypoint / Waypoint).
$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = "object.WAYPOINT"
$ws.Range("B44").Value = "Waypoint"
$ws.Range("C44").Value = "Waypoint"
# No Ukrainian translation supplied for this row yet - make sure column D
# stays genuinely empty (no stray styled-but-blank cell) like the other
# rows that lack a D value (e.g. rows 75/76 before the insert).
$ws.Range("D44").Clear()

# Restore the view/selection state recorded in the saved workbook.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D46").Select() | Out-Null
